# The deck currently has zero slides (no <p:sldIdLst> in presentation.xml).
# Add a new Title Slide (layout 1 == "Title Slide", which exposes the
# ctrTitle / subTitle placeholders) as the presentation's first slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Add(1, 1)
